$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.967.46"
$ws.Range("E2").Value = "  -3.08%  "
$ws.Range("D3").Value = "2.913.86"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'586.68"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "'146.52"
$ws.Range("E6").Value = "  -3.75%  "
$ws.Range("D8").Value = "'0.503"
$ws.Range("E8").Value = "  -3.36%  "
$ws.Range("D9").Value = "2.913.33"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("E10").Value = "  +6.61%  "
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").Value = "'33.72"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.396.54"
$ws.Range("E16").Value = "  -3.95%  "
$ws.Range("D17").Value = "60.918.86"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").Value = "2.912.28"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("D20").Value = "'427.54"
$ws.Range("E20").Value = "  -5.96%  "
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("D22").Value = "'0.673"
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("E23").Value = "  -4.86%  "
$ws.Range("D24").Value = "'80.31"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").Value = "'10.99"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").Value = "'11.89"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").Value = "'2.17"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("E38").Value = "  -3.85%  "
$ws.Range("D39").Value = "'49.36"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("E40").Value = "  -3.88%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  -4.83%  "
$ws.Range("D43").Value = "'0.292"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "'41.80"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").Value = "'377.62"
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "2.673.08"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D50").Value = "'24.88"
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("E51").Value = "  -1.50%  "
